$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update locatorType for rows 3 and 4 to "div"
$ws.Range("E3").Value = "div"
$ws.Range("E4").Value = "div"

# Enable rows 11-14 (set to "Yes")
$ws.Range("B11").Value = "Yes"
$ws.Range("B12").Value = "Yes"
$ws.Range("B13").Value = "Yes"
$ws.Range("B14").Value = "Yes"

# Move the active selection up one row, keeping the same height
$ws.Range("B10:B14").Select()
